$d = $word.ActiveDocument

# 1. Title "Lab 2" -> "Lab 3"
$null = $d.Content.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "3", 1)

# 2. stopwatch_heading
$rng = $d.Content
$found = $rng.Find.Execute("StopWatch version 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: stopwatch_heading" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0C95FC92" w14:textId="77777777" w:rsidR="003F7FCA" w:rsidRPr="00557B06" w:rsidRDefault="003F7FCA" w:rsidP="006940CD"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:i/><w:iCs/></w:rPr><w:t>StopWatch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> version 2</w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 3. completed_stopwatch
$rng = $d.Content
$found = $rng.Find.Execute("Completed version 2 of StopWatch? ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: completed_stopwatch" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="10325CB1" w14:textId="373FB628" w:rsidR="003F7FCA" w:rsidRDefault="00D47471" w:rsidP="00D47471"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Completed version 2 of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StopWatch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">? </w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 4. handlers_merge
$rng = $d.Content
$found = $rng.Find.Execute("Used an arrow function or an anonymous function for the handlers?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: handlers_merge" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4301516C" w14:textId="0170CE79" w:rsidR="00500928" w:rsidRDefault="00500928" w:rsidP="00500928"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>Used an arrow function or an anonymous function for the handlers?</w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 5. created_class_stopwatch
$rng = $d.Content
$found = $rng.Find.Execute("Created a class called StopWatch?  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: created_class_stopwatch" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="03DD678E" w14:textId="77777777" w:rsidR="00500928" w:rsidRDefault="00500928" w:rsidP="00500928"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Created a class called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StopWatch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">?  </w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 6. isfaceup_ismatched
$rng = $d.Content
$found = $rng.Find.Execute(" isFaceUp and isMatched properties", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: isfaceup_ismatched" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6911EB4B" w14:textId="3321A814" w:rsidR="00746974" w:rsidRDefault="00746974" w:rsidP="0011059A"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="00520B5E"><w:t>Has</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>isFaceUp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>isMatched</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> properties</w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 7. fillcards_shufflecards
$rng = $d.Content
$found = $rng.Find.Execute(" a constructor that calls fillCards and shuffleCards?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: fillcards_shufflecards" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7FBB0CBE" w14:textId="46BD4CF8" w:rsidR="005249E2" w:rsidRDefault="00746974" w:rsidP="0011059A"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="00520B5E"><w:t>Has</w:t></w:r><w:r><w:t xml:space="preserve"> a constructor that calls </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fillCards</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shuffleCards</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 8. checkcards
$rng = $d.Content
$found = $rng.Find.Execute(" the methods named above and checkCards", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: checkcards" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4944A7E6" w14:textId="5E48A39A" w:rsidR="00746974" w:rsidRDefault="00746974" w:rsidP="0011059A"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="00520B5E"><w:t>Has</w:t></w:r><w:r><w:t xml:space="preserve"> the methods named above and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>checkCards</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$rng.InsertXML($xml)

# 9. firstpick_secondpick
$rng = $d.Content
$found = $rng.Find.Execute(" properties for matches, tries, firstPick and secondPick?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: firstpick_secondpick" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="63E8E2DA" w14:textId="7CA37CC3" w:rsidR="00746974" w:rsidRDefault="00746974" w:rsidP="0011059A"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="00520B5E"><w:t>Has</w:t></w:r><w:r><w:t xml:space="preserve"> properties for matches, tries, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>firstPick</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>secondPick</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 10. showcards_handleclick_turncardsback
$rng = $d.Content
$found = $rng.Find.Execute("  Has showCards, handleClick and turnCardsBack methods?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: showcards_handleclick_turncardsback" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="59ABB14F" w14:textId="56F25FE5" w:rsidR="005249E2" w:rsidRDefault="00520B5E" w:rsidP="00746974"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">  Has </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>showCards</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>handleClick</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>turnCardsBack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> methods?</w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 11. moved_functionality_init
$rng = $d.Content
$found = $rng.Find.Execute("Moved the functionality of init into the constructor?  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: moved_functionality_init" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0FCD9B7A" w14:textId="77777777" w:rsidR="00441329" w:rsidRDefault="00441329" w:rsidP="00EE05C3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Moved the functionality of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>init</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> into the constructor?  </w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 12. url_github_repo
$rng = $d.Content
$found = $rng.Find.Execute("Url for ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "NOT FOUND: url_github_repo" }
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="763A4223" w14:textId="5841F11D" w:rsidR="00441329" w:rsidRDefault="00441329" w:rsidP="008373F6"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for </w:t></w:r><w:r w:rsidR="002E09A4"><w:t>lab 2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> repo:</w:t></w:r></w:p>'
$rng.InsertXML($xml)

Write-Output "All edits applied"
